$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.005.05"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "'2.055.14"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'248.61"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "'0.659"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'54.57"
$ws.Range("E8").Value = "  +13.98%  "
$ws.Range("D9").Value = "'61.06"
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "'0.378"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  +5.48%  "
$ws.Range("E12").Value = "  +5.67%  "
$ws.Range("D13").Value = "'15.02"
$ws.Range("E13").Value = "  +4.02%  "
$ws.Range("D14").Value = "'2.357.10"
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").Value = "'0.812"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").Value = "'5.21"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "'2.058.27"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "'36.932.42"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "'0.0₃0944"
$ws.Range("E19").Value = "  +13.34%  "
$ws.Range("D20").Value = "'72.27"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "'14.16"
$ws.Range("E21").Value = "  +6.50%  "
$ws.Range("D22").Value = "'5.31"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("D23").Value = "'236.24"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").Value = "'169.81"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").Value = "'9.01"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").Value = "'19.96"
$ws.Range("E28").Value = "  -6.81%  "
$ws.Range("D29").Value = "'1.97"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "'4.53"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").Value = "'0.0620"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").Value = "'1.03"
$ws.Range("E33").Value = "  +9.22%  "
$ws.Range("D34").Value = "'4.35"
$ws.Range("E34").Value = "  +7.02%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.0852"
$ws.Range("E36").Value = "  -6.71%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.26"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -5.84%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("E40").Value = "  +23.38%  "
$ws.Range("D41").Value = "'18.04"
$ws.Range("E41").Value = "  +12.74%  "
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("D44").Value = "'95.99"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "'4.15"
$ws.Range("E46").Value = "  +50.44%  "
$ws.Range("D47").Value = "'2.40"
$ws.Range("E47").Value = "  +7.12%  "
$ws.Range("D48").Value = "'1.293.73"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").Value = "'12.78"
$ws.Range("E50").Value = "  -55.08%  "
$ws.Range("D51").Value = "'6.82"
$ws.Range("E51").Value = "  -2.89%  "
